$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$cell = $table.Cell(4, 4)
$range = $cell.Range
$range.Find.Execute("Fail:401", $true, $false, $false, $false, $false, $true, 0, $false, "Fail:403", 1)
